$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns C and D (Target Suhu, Label Jumlah Awan) - Prediksi Suhu shifts from D to B
$ws.Range("C1:D1").EntireColumn.Delete()
$ws.Range("B1").Value = "Prediksi Suhu"

# Update the remaining values in column B (Prediksi Suhu) with the new predictions
$ws.Range("B2").Value = 29.71184433
$ws.Range("B3").Value = 29.71182699
$ws.Range("B4").Value = 29.71167492
$ws.Range("B5").Value = 29.71165378
$ws.Range("B6").Value = 29.71180773
$ws.Range("B7").Value = 29.71184625
$ws.Range("B8").Value = 29.71176921
$ws.Range("B9").Value = 29.71184452
$ws.Range("B10").Value = 29.71178847
